$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they keep their original
# string representation (matching the source inlineStr cell type) instead of
# being auto-converted to numbers by Excel.
$textCells = @("D5", "D6", "D8", "D12", "D13", "D14", "D16", "D19", "D22", "D23", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D44", "D46", "D49")
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

$ws.Range("D2").Value = '70.038.14'
$ws.Range("E2").Value = '  -1.56%  '
$ws.Range("D3").Value = '3.572.55'
$ws.Range("E3").Value = '  -2.57%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '575.92'
$ws.Range("E5").Value = '  -3.34%  '
$ws.Range("D6").Value = '186.04'
$ws.Range("E6").Value = '  -4.77%  '
$ws.Range("D7").Value = '3.567.49'
$ws.Range("E7").Value = '  -2.59%  '
$ws.Range("D8").Value = '0.620'
$ws.Range("E8").Value = '  -4.38%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("E11").Value = '  -4.13%  '
$ws.Range("D12").Value = '55.04'
$ws.Range("E12").Value = '  -5.69%  '
$ws.Range("D13").Value = '0.0000304'
$ws.Range("E13").Value = '  +2.48%  '
$ws.Range("D14").Value = '9.52'
$ws.Range("E14").Value = '  -4.88%  '
$ws.Range("D15").Value = '4.142.60'
$ws.Range("E15").Value = '  -2.55%  '
$ws.Range("D16").Value = '19.60'
$ws.Range("E16").Value = '  -3.27%  '
$ws.Range("D17").Value = '3.566.01'
$ws.Range("E17").Value = '  -2.71%  '
$ws.Range("D18").Value = '69.918.34'
$ws.Range("E18").Value = '  -1.76%  '
$ws.Range("D19").Value = '12.58'
$ws.Range("E19").Value = '  -2.02%  '
$ws.Range("E21").Value = '  -3.76%  '
$ws.Range("D22").Value = '500.94'
$ws.Range("E22").Value = '  +2.09%  '
$ws.Range("D23").Value = '19.05'
$ws.Range("E24").Value = '  -6.81%  '
$ws.Range("E25").Value = '  -3.18%  '
$ws.Range("D26").Value = '95.38'
$ws.Range("E26").Value = '  +4.35%  '
$ws.Range("D27").Value = '11.30'
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("E28").Value = '  -7.27%  '
$ws.Range("D29").Value = '9.28'
$ws.Range("E29").Value = '  -3.91%  '
$ws.Range("D30").Value = '31.63'
$ws.Range("E30").Value = '  -4.19%  '
$ws.Range("D31").Value = '7.59'
$ws.Range("E31").Value = '  -3.09%  '
$ws.Range("D32").Value = '66.93'
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("D33").Value = '12.05'
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("E34").Value = '  -6.38%  '
$ws.Range("D35").Value = '569.02'
$ws.Range("E35").Value = '  -9.41%  '
$ws.Range("D36").Value = '3.16'
$ws.Range("E36").Value = '  +10.88%  '
$ws.Range("D37").Value = '38.64'
$ws.Range("E37").Value = '  -4.20%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  -5.71%  '
$ws.Range("D40").Value = '0.395'
$ws.Range("E40").Value = '  -4.65%  '
$ws.Range("D41").Value = '3.51'
$ws.Range("E41").Value = '  -2.54%  '
$ws.Range("D42").Value = '3.19'
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("E43").Value = '  -9.93%  '
$ws.Range("D44").Value = '3.00'
$ws.Range("E44").Value = '  -5.73%  '
$ws.Range("D45").Value = '3.231.02'
$ws.Range("E45").Value = '  -3.34%  '
$ws.Range("D46").Value = '3.49'
$ws.Range("E46").Value = '  +4.21%  '
$ws.Range("E47").Value = '  -3.77%  '
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").Value = '0.136'
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("E51").Value = '  -3.65%  '
